$wb = $excel.ActiveWorkbook

# Sheet "展览": update 想去人数 (F column) values for rows 3 and 4
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 146
$wsExhibit.Range("F4").Value = 699

# Sheet "全部类型": update 想去人数 (F column) values for rows 4 and 5
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 146
$wsAll.Range("F5").Value = 699
